$wb = $excel.ActiveWorkbook

# Fold_1
$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 7.22995
$ws.Range("C2").Value = 0.7152999999999999
$ws.Range("D2").Value = 71.01964999999998
$ws.Range("E2").Value = 0.0363
$ws.Range("F2").Value = 27.7257
$ws.Range("G2").Value = 27.7257
$ws.Range("K2").Value = 60.9236
$ws.Range("L2").Value = 27.726
$ws.Range("M2").Value = 33.1976
$ws.Range("N2").Value = 0.191
$ws.Range("O2").Value = 33.00680000000001
$ws.Range("B3").Value = 12.495
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 112.513
$ws.Range("F3").Value = 32.18
$ws.Range("G3").Value = 30.432
$ws.Range("H3").Value = 1.748
$ws.Range("I3").Value = 1.748
$ws.Range("K3").Value = 30.432
$ws.Range("L3").Value = 30.432
$ws.Range("B4").Value = 20.702
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 112.513
$ws.Range("F4").Value = 33.036
$ws.Range("G4").Value = 33.03649863013699
$ws.Range("K4").Value = 33.036
$ws.Range("L4").Value = 33.036
$ws.Range("M4").Value = 0
$ws.Range("O4").Value = 0

# Fold_2
$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 6.946099999999999
$ws.Range("C2").Value = 0.70175
$ws.Range("D2").Value = 71.37609999999998
$ws.Range("E2").Value = 0.08334999999999999
$ws.Range("F2").Value = 27.64249999999999
$ws.Range("G2").Value = 27.64249999999999
$ws.Range("K2").Value = 94.3614
$ws.Range("L2").Value = 27.642
$ws.Range("M2").Value = 66.71939999999999
$ws.Range("N2").Value = 2.934
$ws.Range("O2").Value = 63.7852
$ws.Range("B3").Value = 12.495
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 112.513
$ws.Range("F3").Value = 32.18
$ws.Range("G3").Value = 30.432
$ws.Range("H3").Value = 1.748
$ws.Range("I3").Value = 1.748
$ws.Range("K3").Value = 30.432
$ws.Range("L3").Value = 30.432
$ws.Range("B4").Value = 20.702
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 112.513
$ws.Range("F4").Value = 33.036
$ws.Range("G4").Value = 33.03649863013699
$ws.Range("K4").Value = 33.036
$ws.Range("L4").Value = 33.036
$ws.Range("M4").Value = 0
$ws.Range("O4").Value = 0

# Fold_3
$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 6.944749999999999
$ws.Range("C2").Value = 0.8144
$ws.Range("D2").Value = 64.62134999999999
$ws.Range("E2").Value = 0.06795
$ws.Range("F2").Value = 27.51475
$ws.Range("G2").Value = 27.51475
$ws.Range("K2").Value = 41.2306
$ws.Range("L2").Value = 27.515
$ws.Range("M2").Value = 13.7156
$ws.Range("N2").Value = 10.4554
$ws.Range("O2").Value = 3.2602
$ws.Range("B3").Value = 12.495
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 83.238
$ws.Range("F3").Value = 31.578
$ws.Range("G3").Value = 29.83
$ws.Range("H3").Value = 1.748
$ws.Range("I3").Value = 1.748
$ws.Range("K3").Value = 33.0708
$ws.Range("L3").Value = 29.82999999999999
$ws.Range("B4").Value = 20.702
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 83.238
$ws.Range("F4").Value = 32.435
$ws.Range("G4").Value = 32.43454794520548
$ws.Range("K4").Value = 35.6758
$ws.Range("L4").Value = 32.435
$ws.Range("M4").Value = 3.2408
$ws.Range("N4").Value = 3.2408

# Fold_4
$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 6.124000000000001
$ws.Range("C2").Value = 0.8122000000000001
$ws.Range("D2").Value = 68.59699999999999
$ws.Range("E2").Value = 0.09379999999999999
$ws.Range("F2").Value = 27.3364
$ws.Range("G2").Value = 27.3364
$ws.Range("K2").Value = 81.59760000000001
$ws.Range("L2").Value = 27.336
$ws.Range("M2").Value = 54.2616
$ws.Range("N2").Value = 18.0342
$ws.Range("O2").Value = 36.2274
$ws.Range("B3").Value = 12.495
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 112.513
$ws.Range("F3").Value = 30.432
$ws.Range("G3").Value = 30.432
$ws.Range("K3").Value = 37.4218
$ws.Range("L3").Value = 30.432
$ws.Range("M3").Value = 6.9898
$ws.Range("N3").Value = 6.9898
$ws.Range("B4").Value = 12.495
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 112.513
$ws.Range("F4").Value = 30.432
$ws.Range("G4").Value = 30.43221369863014
$ws.Range("K4").Value = 37.4218
$ws.Range("L4").Value = 30.432
$ws.Range("M4").Value = 6.9898
$ws.Range("N4").Value = 6.9898

# Fold_5
$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 6.6274
$ws.Range("C2").Value = 0.6949500000000002
$ws.Range("D2").Value = 67.77709999999999
$ws.Range("E2").Value = 0.09379999999999999
$ws.Range("F2").Value = 27.46685
$ws.Range("G2").Value = 27.46685
$ws.Range("K2").Value = 135.4388
$ws.Range("L2").Value = 27.467
$ws.Range("M2").Value = 107.9718
$ws.Range("N2").Value = 10.5994
$ws.Range("O2").Value = 97.3724
$ws.Range("B3").Value = 10.141
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 112.513
$ws.Range("F3").Value = 31.934
$ws.Range("G3").Value = 29.685
$ws.Range("H3").Value = 2.249
$ws.Range("I3").Value = 2.249
$ws.Range("K3").Value = 34.1216
$ws.Range("L3").Value = 29.685
$ws.Range("M3").Value = 4.4366
$ws.Range("N3").Value = 3.618
$ws.Range("O3").Value = 0.8186
$ws.Range("B4").Value = 20.702
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 112.513
$ws.Range("F4").Value = 33.036
$ws.Range("G4").Value = 33.03649863013699
$ws.Range("K4").Value = 33.036
$ws.Range("L4").Value = 33.036
